$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1732673267326733
$ws.Range("C2").Value = 0.6138613861386139
$ws.Range("J2").Value = 0.02475247524752475
$ws.Range("P2").Value = 0.1287128712871287
$ws.Range("S2").Value = 0.0594059405940594
$ws.Range("B3").Value = 0.01515151515151515
$ws.Range("C3").Value = 0.06060606060606061
$ws.Range("J3").Value = 0.02272727272727273
$ws.Range("P3").Value = 0.7196969696969697
$ws.Range("S3").Value = 0.1818181818181818
$ws.Range("J4").Value = 0.05405405405405406
$ws.Range("P4").Value = 0.7297297297297297
$ws.Range("S4").Value = 0.2162162162162162
$ws.Range("P5").Value = 1
$ws.Range("B6").Value = 0.04166666666666666
$ws.Range("D6").Value = 0.02083333333333333
$ws.Range("F6").Value = 0.0625
$ws.Range("J6").Value = 0.234375
$ws.Range("O6").Value = 0.005208333333333333
$ws.Range("Q6").Value = 0.1822916666666667
$ws.Range("R6").Value = 0.06770833333333333
$ws.Range("S6").Value = 0.3854166666666667
$ws.Range("B7").Value = 0.07792207792207792
$ws.Range("D7").Value = 0.01948051948051948
$ws.Range("F7").Value = 0.06493506493506493
$ws.Range("J7").Value = 0.07142857142857142
$ws.Range("O7").Value = 0.03246753246753246
$ws.Range("Q7").Value = 0.1233766233766234
$ws.Range("R7").Value = 0.07142857142857142
$ws.Range("S7").Value = 0.538961038961039
$ws.Range("B8").Value = 0.06222222222222222
$ws.Range("D8").Value = 0.01333333333333333
$ws.Range("F8").Value = 0.07333333333333333
$ws.Range("J8").Value = 0.1088888888888889
$ws.Range("O8").Value = 0.02
$ws.Range("Q8").Value = 0.1422222222222222
$ws.Range("R8").Value = 0.1022222222222222
$ws.Range("S8").Value = 0.4777777777777778
$ws.Range("B9").Value = 0.06976744186046512
$ws.Range("D9").Value = 0.005813953488372093
$ws.Range("F9").Value = 0.0755813953488372
$ws.Range("J9").Value = 0.09883720930232558
$ws.Range("O9").Value = 0.005813953488372093
$ws.Range("Q9").Value = 0.1569767441860465
$ws.Range("R9").Value = 0.09302325581395349
$ws.Range("S9").Value = 0.4941860465116279
$ws.Range("B10").Value = 0.09872922776148582
$ws.Range("D10").Value = 0.02248289345063539
$ws.Range("E10").Value = 0.002932551319648094
$ws.Range("F10").Value = 0.06451612903225806
$ws.Range("J10").Value = 0.09481915933528837
$ws.Range("O10").Value = 0.02346041055718475
$ws.Range("Q10").Value = 0.1994134897360704
$ws.Range("R10").Value = 0.09775171065493646
$ws.Range("S10").Value = 0.3958944281524927
$ws.Range("F11").Value = 0.004273504273504274
$ws.Range("G11").Value = 0.1367521367521368
$ws.Range("J11").Value = 0.07692307692307693
$ws.Range("K11").Value = 0.188034188034188
$ws.Range("L11").Value = 0.5811965811965812
$ws.Range("S11").Value = 0.01282051282051282
$ws.Range("J12").Value = 0.1956521739130435
$ws.Range("K12").Value = 0.01449275362318841
$ws.Range("L12").Value = 0.02898550724637681
$ws.Range("S12").Value = 0.04347826086956522
$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.2444444444444444
$ws.Range("S13").Value = 0.08888888888888889
$ws.Range("F15").Value = 0.01507537688442211
$ws.Range("H15").Value = 0.1507537688442211
$ws.Range("I15").Value = 0.07035175879396985
$ws.Range("J15").Value = 0.3819095477386935
$ws.Range("K15").Value = 0.05527638190954774
$ws.Range("M15").Value = 0.005025125628140704
$ws.Range("O15").Value = 0.05527638190954774
$ws.Range("S15").Value = 0.2663316582914573
$ws.Range("F16").Value = 0.00684931506849315
$ws.Range("H16").Value = 0.2054794520547945
$ws.Range("I16").Value = 0.07534246575342465
$ws.Range("J16").Value = 0.3424657534246575
$ws.Range("K16").Value = 0.1164383561643836
$ws.Range("M16").Value = 0.03424657534246575
$ws.Range("O16").Value = 0.04794520547945205
$ws.Range("S16").Value = 0.1712328767123288
$ws.Range("F17").Value = 0.01123595505617977
$ws.Range("H17").Value = 0.1938202247191011
$ws.Range("I17").Value = 0.1151685393258427
$ws.Range("J17").Value = 0.4044943820224719
$ws.Range("K17").Value = 0.06741573033707865
$ws.Range("M17").Value = 0.01966292134831461
$ws.Range("O17").Value = 0.07303370786516854
$ws.Range("S17").Value = 0.1151685393258427
$ws.Range("F18").Value = 0.03174603174603174
$ws.Range("H18").Value = 0.2433862433862434
$ws.Range("I18").Value = 0.07407407407407407
$ws.Range("J18").Value = 0.4021164021164021
$ws.Range("K18").Value = 0.1005291005291005
$ws.Range("M18").Value = 0.02116402116402116
$ws.Range("N18").Value = 0.005291005291005291
$ws.Range("O18").Value = 0.03174603174603174
$ws.Range("S18").Value = 0.08994708994708994
$ws.Range("F19").Value = 0.01877133105802048
$ws.Range("H19").Value = 0.2414675767918089
$ws.Range("I19").Value = 0.07935153583617748
$ws.Range("J19").Value = 0.3506825938566553
$ws.Range("K19").Value = 0.09385665529010238
$ws.Range("M19").Value = 0.02815699658703072
$ws.Range("N19").Value = 0.0008532423208191126
$ws.Range("O19").Value = 0.07081911262798635
$ws.Range("S19").Value = 0.1160409556313993
